$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price values so Excel does not
# coerce them to numbers (which would lose significant trailing zeros,
# e.g. "1.0000" -> 1). These cells are plain inline strings in the sheet.
$numericTextCells = @("D5","D6","D7","D8","D9","D10","D11","D12","D13","D15","D17","D19","D21","D22","D24","D25","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "30.134.26"
$ws.Range("E2").Value = "  +5.89%  "
$ws.Range("D3").Value = "1.921.98"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  -0.93%  "
$ws.Range("D5").Value = "327.22"
$ws.Range("E5").Value = "  +3.60%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "0.5176"
$ws.Range("E7").Value = "  +1.94%  "
$ws.Range("D8").Value = "0.4050"
$ws.Range("E8").Value = "  +3.97%  "
$ws.Range("D9").Value = "0.08478"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").Value = "1.128"
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("D11").Value = "42.86"
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("D12").Value = "22.19"
$ws.Range("E12").Value = "  +8.88%  "
$ws.Range("D13").Value = "6.364"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "1.926.54"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").Value = "7.388"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "96.02"
$ws.Range("E17").Value = "  +5.26%  "
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").Value = "0.06743"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("E20").Value = "  +3.45%  "
$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "6.065"
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").Value = "30.128.59"
$ws.Range("E23").Value = "  +5.74%  "
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").Value = "2.198"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("D26").Value = "2.147.11"
$ws.Range("E26").Value = "  +2.93%  "
$ws.Range("D27").Value = "21.29"
$ws.Range("E27").Value = "  +3.43%  "
$ws.Range("D28").Value = "160.76"
$ws.Range("D29").Value = "2.464"
$ws.Range("E29").Value = "  +4.19%  "
$ws.Range("D30").Value = "129.34"
$ws.Range("E30").Value = "  +2.94%  "
$ws.Range("D31").Value = "1.083"
$ws.Range("E31").Value = "  +4.84%  "
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").Value = "6.099"
$ws.Range("E33").Value = "  +5.99%  "
$ws.Range("D34").Value = "3.656"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("D35").Value = "0.02518"
$ws.Range("E35").Value = "  +2.55%  "
$ws.Range("D36").Value = "0.06624"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").Value = "0.2217"
$ws.Range("E37").Value = "  +2.78%  "
$ws.Range("D38").Value = "1.242"
$ws.Range("E38").Value = "  +4.85%  "
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("D40").Value = "9.035"
$ws.Range("E40").Value = "  +2.41%  "
$ws.Range("D41").Value = "0.6580"
$ws.Range("E41").Value = "  +3.07%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.251"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "11.64"
$ws.Range("E43").Value = "  +5.09%  "
$ws.Range("D44").Value = "0.6182"
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").Value = "13.28"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("D46").Value = "3.761"
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("D47").Value = "2.071"
$ws.Range("E47").Value = "  +3.59%  "
$ws.Range("D48").Value = "1.247"
$ws.Range("E48").Value = "  +2.78%  "
$ws.Range("D49").Value = "125.98"
$ws.Range("E49").Value = "  +3.69%  "
$ws.Range("D50").Value = "1.164"
$ws.Range("E50").Value = "  +3.28%  "
$ws.Range("D51").Value = "79.67"
$ws.Range("E51").Value = "  +4.63%  "
